$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approve")

# Change the endpoint text for the first block (row 3) to the new API path
$ws.Range("E3").Value = "/conclusions/approve1"

# Clear out the second block (rows 5-7), which previously held "Assert400" data
$ws.Range("A5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""

# Update selection to match the new active area
$ws.Range("A5:E7").Select()
